# Apply the "adding more fields and updating DB structure" change:
#  - Insert a new "hash_id" field row into the OS_Directory_Item table (column F),
#    right after item_id / before mode, shifting mode/size/is_link/is_hidden down one row.
#  - Append a new set of rows describing the fields of a Windows PE "file info" record
#    (comments, company_name, ... product_version) below the OS_Directory_Item table.
#  - Add a new "Hash" table in column J with fields id, md5, sha1, sha256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing OS_Directory_Item fields (mode, size, is_link, is_hidden) down by one
# row to make room for the new "hash_id" field, without doing a full-row insert (which would
# also drag column A's background styling down with it).
$ws.Range("F10").Value = $ws.Range("F9").Value2
$ws.Range("F9").Value = $ws.Range("F8").Value2
$ws.Range("F8").Value = $ws.Range("F7").Value2
$ws.Range("F7").Value = $ws.Range("F6").Value2

# --- New "Hash" table header (column J), formats copied from the "Item_ADS" column (I) ---
$ws.Range("I1:I5").Copy()
$ws.Range("J1:J5").PasteSpecial(-4122)  # xlPasteFormats

# The order cell values are written below matters: new text values are appended to the
# workbook's shared-string table in first-use order, and must match the target file.

# Extra fields appended to the OS_Directory_Item table (rows 11-16, 18-29)
$ws.Range("F11").Value = "comments"
$ws.Range("F12").Value = "company_name"
$ws.Range("F13").Value = "file_build_part"
$ws.Range("F14").Value = "file_description"
$ws.Range("F15").Value = "file_major_part"
$ws.Range("F16").Value = "file_minor_part"
$ws.Range("F18").Value = "file_private_part"
$ws.Range("F19").Value = "file_version"
$ws.Range("F20").Value = "internal_name"
$ws.Range("F21").Value = "is_debug"
$ws.Range("F22").Value = "is_patched"
$ws.Range("F23").Value = "is_private_build"
$ws.Range("F24").Value = "is_prerelease"
$ws.Range("F25").Value = "is_special_build"
$ws.Range("F26").Value = "language"
$ws.Range("F27").Value = "legal_copyright"
$ws.Range("F28").Value = "legal_trademarks"
$ws.Range("F29").Value = "original_filename"
$ws.Range("F17").Value = "filename"
$ws.Range("F30").Value = "private_build"
$ws.Range("F31").Value = "product_build_part"
$ws.Range("F32").Value = "product_major_part"
$ws.Range("F33").Value = "product_minor_part"
$ws.Range("F34").Value = "product_name"
$ws.Range("F35").Value = "product_private_part"
$ws.Range("F36").Value = "product_version"

# New "Hash" table (column J)
$ws.Range("J3").Value = "md5"
$ws.Range("J4").Value = "sha1"
$ws.Range("J5").Value = "sha256"
$ws.Range("J1").Value = "Hash"

# New field on OS_Directory_Item referencing the Hash table
$ws.Range("F6").Value = "hash_id"

# id field for the Hash table (reuses the already-existing "id" shared string)
$ws.Range("J2").Value = "id"

# Update the selected/active cell to match the authored workbook
$ws.Range("F7").Select()
